# "Make problems more interesting": swap in new character names/scenarios
# for the question line and the "Answer" line, and widen the blank lines
# for "Number sentence" and "Answer" in each of the 7 arithmetic word
# problems in the document.
#
# Each problem paragraph holds three lines (question, "Number sentence:
# ___", "Answer: ___") separated by manual line breaks (<w:br/>), all
# inside one paragraph. Find/Replace represents a <w:br/> as a vertical-tab
# character (chr 11) in Range.Text, so each paragraph's old/new text is
# built by joining the three old/new lines with that character and doing
# one Find.Execute per paragraph, scoped to that paragraph's Range.

$d = $word.ActiveDocument
$BR = [char]11

function Replace-Problem($paraIndex, $oldQuestion, $newQuestion, $oldAnswer, $newAnswer) {
    $oldNumberSentence = "Number sentence: __________________________________________________"
    $newNumberSentence = "Number sentence: ______________________________________________________________________"

    $old = $oldQuestion + $BR + $oldNumberSentence + $BR + $oldAnswer
    $new = $newQuestion + $BR + $newNumberSentence + $BR + $newAnswer

    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.Find.Execute($old, $true, $true, $false, $false, $false, `
                     $true, 1, $false, $new, 2)
}

# Problem 1 (paragraph 2): Lily -> Barnaby
Replace-Problem 2 `
    "Lily has four crayons and gets three more. How many crayons does Lily have in total?" `
    "Barnaby the bear found seven honey pots. Gertrude the goose gave him two more. How many honey pots does Barnaby have in all?" `
    "Answer: Lily has ____________ crayons in total." `
    "Answer: Barnaby has ________________________________________ honey pots in all."

# Problem 2 (paragraph 3): David -> Penelope
Replace-Problem 3 `
    "David had eleven toy cars, but he lost one of them. How many toy cars does David have now?" `
    "Penelope the penguin had eleven shiny pebbles. She lost three pebbles down a crack. How many pebbles does Penelope have now?" `
    "Answer: David now has ____________ toy cars." `
    "Answer: Penelope now has ________________________________________ pebbles."

# Problem 3 (paragraph 4): Sarah -> Rupert
Replace-Problem 4 `
    "Sarah baked nine cookies, and her friend ate four of them. How many cookies are left?" `
    "Rupert the rabbit collected eight carrots from the garden. He ate four of them for lunch. How many carrots does Rupert have left?" `
    "Answer: There are ____________ cookies left." `
    "Answer: Rupert has ________________________________________ carrots left."

# Problem 4 (paragraph 5): apples -> Flora
Replace-Problem 5 `
    "There are eight apples on a tree. Four more apples grow. How many apples are on the tree?" `
    "Flora the fairy had nine sparkly wands. She found six more wands under a mushroom. How many wands does Flora have altogether?" `
    "Answer: There are now ____________ apples on the tree." `
    "Answer: Flora has ________________________________________ wands altogether."

# Problem 5 (paragraph 6): Ben -> Zorp
Replace-Problem 6 `
    "Ben has seven building blocks and his sister gives him two more. How many blocks does Ben have?" `
    "Zorp the alien had twelve wiggly worms. Three of them wriggled away. How many wiggly worms does Zorp have now?" `
    "Answer: Ben has ____________ blocks." `
    "Answer: Zorp has ________________________________________ wiggly worms now."

# Problem 6 (paragraph 7): Chloe -> Brenda/Cecil
Replace-Problem 7 `
    "Chloe had twelve beads. She used three beads to make a bracelet. How many beads does Chloe have now?" `
    "Brenda the badger baked five pies. Cecil the squirrel baked three pies. How many pies did they bake together?" `
    "Answer: Chloe now has ____________ beads." `
    "Answer: They baked ________________________________________ pies together."

# Problem 7 (paragraph 8): Emily -> Horace
Replace-Problem 8 `
    "Emily has five stickers. She gives one sticker to her friend. How many stickers does Emily have left?" `
    "Horace the hedgehog had fourteen spiky pine cones. He gave five pine cones to his friend. How many pine cones does Horace have now?" `
    "Answer: Emily has ____________ stickers left." `
    "Answer: Horace has ________________________________________ pine cones now."
